$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the indicator metadata cells in column B ---

# B4: Indicator name (wording tweak, no trailing period, no leading/trailing space)
$ws.Range("B4").Value = "9.1.2 Объем пассажирских и грузовых перевозок в разбивке по видам транспорта"

# B6: Organisation name / department
$ws.Range("B6").Value = "Национальный статистический комитет КР`n(Управление статистики торговли и услуг, Информационно-коммуникационных технологий и туризма)"
$ws.Range("B6").Font.Name = "Calibri"

# B7: Contact person (new name, keeps trailing line break) - also becomes wrapped like B6
$ws.Range("B7").Value = "Текеева Л.А.`n"
$ws.Range("B7").Font.Name = "Calibri"
$ws.Range("B7").WrapText = $true

# B8: Contact e-mail
$ws.Range("B8").Value = "Ltekeeva@stat.kg"

# B9: Contact phone number
$ws.Range("B9").Value = "0 (312) 32-47-25"
$ws.Range("B9").Font.Name = "Calibri"

# B10: Organisation website
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Font.Name = "Calibri"

# --- Move the active selection to B7 ---
$null = $ws.Range("B7").Select()
